# Add team record (Wins / Losses / Ties) columns to the roster sheet.
# Mirrors the commit "Added team record to data" — new header cells in
# AD1:AF1 (same header style as the existing headers) and a constant
# 94-67-0 record value written into AD:AF for every data row (2-44).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header labels.
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Give the new headers the same formatting (bold, centered, thin border)
# as the rest of row 1 by copying the format from the adjacent header cell.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Determine the last populated data row (falls back to 44 as seen in the source file).
$lastRow = $ws.UsedRange.Rows.Count
if ($lastRow -lt 44) {
    $lastRow = 44
}

# Fill the team record for every data row.
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 94   # AD - Wins
    $ws.Cells.Item($r, 31).Value = 67   # AE - Losses
    $ws.Cells.Item($r, 32).Value = 0    # AF - Ties
}
